$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vtab = [char]11

$t.Cell(1,1).Range.Text = "84 x 78" + $vtab + "  7    8" + $vtab + "  ----" + $vtab + "8|    |" + $vtab + "4|    |"
$t.Cell(1,2).Range.Text = "93 x 87" + $vtab + "  8    7" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "3|    |"
$t.Cell(1,3).Range.Text = "33 x 43" + $vtab + "  4    3" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "3|    |"
$t.Cell(2,1).Range.Text = "44 x 63" + $vtab + "  6    3" + $vtab + "  ----" + $vtab + "4|    |" + $vtab + "4|    |"
$t.Cell(2,2).Range.Text = "54 x 65" + $vtab + "  6    5" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "4|    |"
$t.Cell(2,3).Range.Text = "22 x 66" + $vtab + "  6    6" + $vtab + "  ----" + $vtab + "2|    |" + $vtab + "2|    |"
$t.Cell(3,1).Range.Text = "92 x 40" + $vtab + "  4    0" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "2|    |"
$t.Cell(3,2).Range.Text = "99 x 96" + $vtab + "  9    6" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "9|    |"
$t.Cell(3,3).Range.Text = "58 x 29" + $vtab + "  2    9" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "8|    |"
$t.Cell(4,1).Range.Text = "91 x 38" + $vtab + "  3    8" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "1|    |"
$t.Cell(4,2).Range.Text = "22 x 42" + $vtab + "  4    2" + $vtab + "  ----" + $vtab + "2|    |" + $vtab + "2|    |"
$t.Cell(4,3).Range.Text = "62 x 87" + $vtab + "  8    7" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "2|    |"
$t.Cell(5,1).Range.Text = "14 x 35" + $vtab + "  3    5" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "4|    |"
$t.Cell(5,2).Range.Text = "93 x 74" + $vtab + "  7    4" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "3|    |"
$t.Cell(5,3).Range.Text = "19 x 72" + $vtab + "  7    2" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "9|    |"
